$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
# aggiornamento fino a 20/09/2021
$data = @(
    @(44449, 0, 5, 57.49770009199631),
    @(44450, 2, 5, 57.49770009199631),
    @(44451, 1, 6, 68.99724011039559),
    @(44452, 0, 5, 57.49770009199631),
    @(44453, 1, 4, 45.99816007359706),
    @(44454, 0, 4, 45.99816007359706),
    @(44455, 1, 5, 57.49770009199631),
    @(44456, 1, 6, 68.99724011039559),
    @(44457, 1, 5, 57.49770009199631),
    @(44458, 1, 5, 57.49770009199631),
    @(44459, 1, 6, 68.99724011039559)
)

$startRow = 375

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Carry the date-formatted style from the last existing row down to the new row (col A).
    $ws.Range("A374").Copy($ws.Range("A$r"))

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
